# Swap the (D,E) column-pair with the (F,G) column-pair for every row
# (header row included) on the active sheet.
#
# Before: A=code B=name C=status D=category-code E=group-code
#         F=group-name  G=category-name
# After:  A=code B=name C=status D=group-name    E=category-name
#         F=group-code  G=category-code

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

$srcRange = $ws.Range("D1:G$lastRow")
$vals = $srcRange.Value()

# Force the destination columns to be stored as text so that numeric-looking
# codes (e.g. "110", "111") round-trip as strings instead of being coerced
# into numbers when the swapped array is written back.
$srcRange.NumberFormat = "@"

$rows = $vals.GetLength(0)
$new = New-Object 'object[,]' $rows,4

for ($i = 1; $i -le $rows; $i++) {
    $new[$i-1,0] = $vals[$i,3]   # new D = old F (group-name / group-code header)
    $new[$i-1,1] = $vals[$i,4]   # new E = old G (category-name / category-code header)
    $new[$i-1,2] = $vals[$i,2]   # new F = old E (group-code / group-name)
    $new[$i-1,3] = $vals[$i,1]   # new G = old D (category-code / category-name)
}

$srcRange.Value = $new
